$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh with new TPM data: rows 2-7, columns E (Ligand-expressing cells)
# through T (Edge total expression derived specificity). Columns A-D
# (Sending cluster / Ligand symbol / Receptor symbol / Target cluster)
# are unchanged.

$newValues = @{
    2 = @{ E=2; F=0.6666666666666666; G=0.2892773333333333; H=0.8678319999999999; I=1; J=1; K=2; L=1; M=3.6833285;          N=7.366657;          O=0.3319404283605227; P=0.2657978481314736; Q=1.065503446270667;  R=6.393020677624;    S=0.3319404283605227; T=0.2657978481314736 }
    3 = @{ E=2; F=0.6666666666666666; G=0.2892773333333333; H=0.8678319999999999; I=1; J=1; K=2; L=0.6666666666666666; M=1.145672333333333;   N=3.437017;           O=0.1032476373170262; P=0.1240117087834133; Q=0.3314170374604445;  R=2.982753337144;    S=0.1032476373170262; T=0.1240117087834133 }
    4 = @{ E=2; F=0.6666666666666666; G=0.2892773333333333; H=0.8678319999999999; I=1; J=1; K=2; L=0.6666666666666666; M=0.5012456666666667;  N=1.503737;           O=0.04517210487937449; P=0.05425664025835297; Q=0.1449990097982222; R=1.304991088184;    S=0.04517210487937449; T=0.05425664025835297 }
    5 = @{ E=2; F=0.6666666666666666; G=0.2892773333333333; H=0.8678319999999999; I=1; J=1; K=2; L=1; M=1.8904715;           N=3.780943;           O=0.1703687084965025; P=0.1364209726756327; Q=0.5468705542626667;  R=3.281223325576;    S=0.1703687084965025; T=0.1364209726756327 }
    6 = @{ E=2; F=0.6666666666666666; G=0.2892773333333333; H=0.8678319999999999; I=1; J=1; K=3; L=1; M=3.028466;            N=9.085398;           O=0.2729244218416247; P=0.3278120914029245; Q=0.8760665685706667;  R=7.884599117135999; S=0.2729244218416247; T=0.3278120914029245 }
    7 = @{ E=2; F=0.6666666666666666; G=0.2892773333333333; H=0.8678319999999999; I=1; J=1; K=3; L=1; M=0.8471700000000001;  N=2.54151;            O=0.07634669910494926; P=0.09170073874820307; Q=0.24506707848;       R=2.20560370632;     S=0.07634669910494926; T=0.09170073874820307 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
